$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 1.82
$ws.Range("I2").Value = 3.75
$ws.Range("J2").Value = 2.38
$ws.Range("N2").Value = 19
$ws.Range("Q2").Value = 1.44
$ws.Range("R2").Value = 2.75
$ws.Range("U2").Value = 1.44
$ws.Range("V2").Value = 2.63
$ws.Range("Z2").Value = 17
$ws.Range("AB2").Value = 17
$ws.Range("AD2").Value = 8.5
$ws.Range("AM2").Value = 23
$ws.Range("AN2").Value = 4.33

# Row 7 updates
$ws.Range("G7").Value = 3.2
$ws.Range("I7").Value = 2.2
$ws.Range("J7").Value = 3.75
$ws.Range("L7").Value = 2.88
$ws.Range("X7").Value = 17
$ws.Range("Y7").Value = 12
$ws.Range("AA7").Value = 26
$ws.Range("AH7").Value = 8
$ws.Range("AJ7").Value = 9
$ws.Range("AL7").Value = 17
$ws.Range("AR7").Value = 81
$ws.Range("AY7").Value = 12
